$wb = $excel.ActiveWorkbook

# ---- Sheet: ALC ----
$ws = $wb.Worksheets.Item("ALC")
# Row 64
$ws.Range("H64").Value = 704730
$ws.Range("I64").Value = 959438.8
$ws.Range("J64").Value = 4280.75
$ws.Range("K64").Value = 959438.8
$ws.Range("L64").Value = 4280.75
$ws.Range("M64").Value = -959190.8
# Row 67
$ws.Range("H67").Value = 704730
$ws.Range("I67").Value = 959438.8
$ws.Range("J67").Value = 4280.75
$ws.Range("K67").Value = 959438.8
$ws.Range("L67").Value = 4280.75
$ws.Range("M67").Value = -958580.8
# Row 86
$ws.Range("H86").Value = 4124.2
$ws.Range("I86").Value = 4000
$ws.Range("J86").Value = 4207
$ws.Range("K86").Value = 4000
$ws.Range("L86").Value = 4207
$ws.Range("M86").Value = -2877
$ws.Range("N86").Value = -6453
# Row 89
$ws.Range("H89").Value = 4124.2
$ws.Range("I89").Value = 4000
$ws.Range("J89").Value = 4207
$ws.Range("K89").Value = 20000
$ws.Range("L89").Value = 21035
$ws.Range("M89").Value = -14384
$ws.Range("N89").Value = -32267
# Row 121
$ws.Range("H121").Value = 778.6667
$ws.Range("I121").Value = 1800
$ws.Range("J121").Value = 718.58826
$ws.Range("K121").Value = 5400
$ws.Range("L121").Value = 2155.76478
$ws.Range("M121").Value = -3653
$ws.Range("N121").Value = -5649.76478
# Row 132
$ws.Range("H132").Value = 2118.75
$ws.Range("I132").Value = 1900.4445
$ws.Range("J132").Value = 2773.6667
$ws.Range("K132").Value = 5701.333500000001
$ws.Range("L132").Value = 8321.000100000001
$ws.Range("M132").Value = -3171.333500000001
$ws.Range("N132").Value = -13381.0001
# Row 135
$ws.Range("H135").Value = 2427.05
$ws.Range("I135").Value = 1252.1818
$ws.Range("J135").Value = 7965.7144
$ws.Range("K135").Value = 11269.6362
$ws.Range("L135").Value = 71691.4296
$ws.Range("M135").Value = -8734.636200000001
$ws.Range("N135").Value = -76761.4296
# Row 137
$ws.Range("H137").Value = 1227.579
$ws.Range("I137").Value = 715.625
$ws.Range("J137").Value = 1599.909
$ws.Range("K137").Value = 2146.875
$ws.Range("L137").Value = 4799.727000000001
$ws.Range("M137").Value = 403.125
$ws.Range("N137").Value = -9899.727000000001

# ---- Sheet: ARM ----
$ws = $wb.Worksheets.Item("ARM")
# Row 32
$ws.Range("H32").Value = 2043.495
$ws.Range("I32").Value = 1770.7174
$ws.Range("J32").Value = 5628.5713
$ws.Range("K32").Value = 1770.7174
$ws.Range("L32").Value = 5628.5713
$ws.Range("M32").Value = -1483.7174
# Row 61
$ws.Range("H61").Value = 2253.8572
$ws.Range("I61").Value = 2289
$ws.Range("J61").Value = 2199.5454
$ws.Range("K61").Value = 2289
$ws.Range("L61").Value = 2199.5454
$ws.Range("M61").Value = -2077
$ws.Range("N61").Value = -2623.5454
# Row 74
$ws.Range("H74").Value = 3052.0588
$ws.Range("I74").Value = 3480.7693
$ws.Range("J74").Value = 1658.75
$ws.Range("K74").Value = 3480.7693
$ws.Range("L74").Value = 1658.75
$ws.Range("M74").Value = -2606.7693
$ws.Range("N74").Value = -3406.75
# Row 77
$ws.Range("H77").Value = 3052.0588
$ws.Range("I77").Value = 3480.7693
$ws.Range("J77").Value = 1658.75
$ws.Range("K77").Value = 17403.8465
$ws.Range("L77").Value = 8293.75
$ws.Range("M77").Value = -13035.8465
$ws.Range("N77").Value = -17029.75
# Row 133
$ws.Range("H133").Value = 68351.836
$ws.Range("I133").Value = 0
$ws.Range("J133").Value = 68351.836
$ws.Range("K133").Value = 0
$ws.Range("L133").Value = 68351.836
$ws.Range("N133").Value = -73411.836
# Row 135
$ws.Range("H135").Value = 0
$ws.Range("I135").Value = 0
$ws.Range("J135").Value = 0
$ws.Range("K135").Value = 0
$ws.Range("L135").Value = 0
$ws.Range("N135").ClearContents()
# Row 136
$ws.Range("H136").Value = 2253.8572
$ws.Range("I136").Value = 2289
$ws.Range("J136").Value = 2199.5454
$ws.Range("K136").Value = 6867
$ws.Range("L136").Value = 6598.6362
$ws.Range("M136").Value = -4317
$ws.Range("N136").Value = -11698.6362
# Row 138
$ws.Range("H138").Value = 49096.668
$ws.Range("I138").Value = 0
$ws.Range("J138").Value = 49096.668
$ws.Range("K138").Value = 0
$ws.Range("L138").Value = 49096.668
$ws.Range("N138").Value = -59376.668
# Row 139
$ws.Range("H139").Value = 30000
$ws.Range("I139").Value = 0
$ws.Range("J139").Value = 30000
$ws.Range("K139").Value = 0
$ws.Range("L139").Value = 30000
$ws.Range("N139").Value = -40280

# ---- Sheet: BSM ----
$ws = $wb.Worksheets.Item("BSM")
# Row 22
$ws.Range("H22").Value = 400
$ws.Range("I22").Value = 400
$ws.Range("J22").Value = 0
$ws.Range("K22").Value = 400
$ws.Range("L22").Value = 0
$ws.Range("M22").Value = -227
$ws.Range("N22").ClearContents()
# Row 69
$ws.Range("H69").Value = 10000
$ws.Range("I69").Value = 10000
$ws.Range("J69").Value = 0
$ws.Range("K69").Value = 10000
$ws.Range("L69").Value = 0
$ws.Range("M69").Value = -9189
# Row 72
$ws.Range("H72").Value = 10000
$ws.Range("I72").Value = 10000
$ws.Range("J72").Value = 0
$ws.Range("K72").Value = 30000
$ws.Range("L72").Value = 0
$ws.Range("M72").Value = -25944

# ---- Sheet: CRP ----
$ws = $wb.Worksheets.Item("CRP")
# Row 31
$ws.Range("H31").Value = 4147.05
$ws.Range("I31").Value = 2715.1333
$ws.Range("J31").Value = 8442.799999999999
$ws.Range("K31").Value = 2715.1333
$ws.Range("L31").Value = 8442.799999999999
$ws.Range("M31").Value = -2420.1333
$ws.Range("N31").Value = -9032.799999999999
# Row 34
$ws.Range("H34").Value = 4147.05
$ws.Range("I34").Value = 2715.1333
$ws.Range("J34").Value = 8442.799999999999
$ws.Range("K34").Value = 2715.1333
$ws.Range("L34").Value = 8442.799999999999
$ws.Range("M34").Value = -2513.1333
$ws.Range("N34").Value = -8846.799999999999
# Row 58
$ws.Range("H58").Value = 1350.9474
$ws.Range("I58").Value = 807.5
$ws.Range("J58").Value = 1746.1818
$ws.Range("K58").Value = 807.5
$ws.Range("L58").Value = 1746.1818
$ws.Range("M58").Value = -604.5
$ws.Range("N58").Value = -2152.1818
# Row 94
$ws.Range("H94").Value = 4760.8237
$ws.Range("I94").Value = 700
$ws.Range("J94").Value = 5014.625
$ws.Range("K94").Value = 700
$ws.Range("L94").Value = 5014.625
$ws.Range("M94").Value = -249
$ws.Range("N94").Value = -5916.625
# Row 136
$ws.Range("H136").Value = 1350.9474
$ws.Range("I136").Value = 807.5
$ws.Range("J136").Value = 1746.1818
$ws.Range("K136").Value = 2422.5
$ws.Range("L136").Value = 5238.5454
$ws.Range("M136").Value = 127.5
$ws.Range("N136").Value = -10338.5454
# Row 140
$ws.Range("H140").Value = 87618.336
$ws.Range("I140").Value = 0
$ws.Range("J140").Value = 87618.336
$ws.Range("K140").Value = 0
$ws.Range("L140").Value = 87618.336
$ws.Range("N140").Value = -97978.336

# ---- Sheet: CUL ----
$ws = $wb.Worksheets.Item("CUL")
# Row 64
$ws.Range("H64").Value = 842.7143
$ws.Range("I64").Value = 428.33334
$ws.Range("J64").Value = 1153.5
$ws.Range("K64").Value = 1285.00002
$ws.Range("L64").Value = 3460.5
$ws.Range("M64").Value = -1015.00002
$ws.Range("N64").Value = -4000.5
# Row 67
$ws.Range("H67").Value = 842.7143
$ws.Range("I67").Value = 428.33334
$ws.Range("J67").Value = 1153.5
$ws.Range("K67").Value = 1285.00002
$ws.Range("L67").Value = 3460.5
$ws.Range("M67").Value = -349.0000199999999
$ws.Range("N67").Value = -5332.5
# Row 74
$ws.Range("H74").Value = 7617.75
$ws.Range("I74").Value = 3500
$ws.Range("J74").Value = 9676.625
$ws.Range("K74").Value = 10500
$ws.Range("L74").Value = 29029.875
$ws.Range("M74").Value = -9439
$ws.Range("N74").Value = -31151.875
# Row 77
$ws.Range("H77").Value = 7617.75
$ws.Range("I77").Value = 3500
$ws.Range("J77").Value = 9676.625
$ws.Range("K77").Value = 31500
$ws.Range("L77").Value = 87089.625
$ws.Range("M77").Value = -26196
$ws.Range("N77").Value = -97697.625
# Row 136
$ws.Range("H136").Value = 1576.7273
$ws.Range("I136").Value = 1140.9375
$ws.Range("J136").Value = 2738.8333
$ws.Range("K136").Value = 3422.8125
$ws.Range("L136").Value = 8216.499899999999
$ws.Range("M136").Value = 1677.1875
# Row 138
$ws.Range("H138").Value = 4186.409
$ws.Range("I138").Value = 5021.5
$ws.Range("J138").Value = 2725
$ws.Range("K138").Value = 15064.5
$ws.Range("L138").Value = 8175
$ws.Range("M138").Value = -9924.5
$ws.Range("N138").Value = -18455

# ---- Sheet: GSM ----
$ws = $wb.Worksheets.Item("GSM")
# Row 140
$ws.Range("H140").Value = 41329.8
$ws.Range("I140").Value = 0
$ws.Range("J140").Value = 41329.8
$ws.Range("K140").Value = 0
$ws.Range("L140").Value = 41329.8
$ws.Range("N140").Value = -51689.8

# ---- Sheet: LTW ----
$ws = $wb.Worksheets.Item("LTW")
# Row 68
$ws.Range("H68").Value = 2167.4644
$ws.Range("I68").Value = 2124.5
$ws.Range("J68").Value = 2274.875
$ws.Range("K68").Value = 2124.5
$ws.Range("L68").Value = 2274.875
$ws.Range("M68").Value = -1375.5
$ws.Range("N68").Value = -3772.875
# Row 71
$ws.Range("H71").Value = 2167.4644
$ws.Range("I71").Value = 2124.5
$ws.Range("J71").Value = 2274.875
$ws.Range("K71").Value = 10622.5
$ws.Range("L71").Value = 11374.375
$ws.Range("M71").Value = -6878.5
$ws.Range("N71").Value = -18862.375
# Row 100
$ws.Range("H100").Value = 3729.8572
$ws.Range("I100").Value = 3498.25
$ws.Range("J100").Value = 3822.5
$ws.Range("K100").Value = 3498.25
$ws.Range("L100").Value = 3822.5
$ws.Range("M100").Value = -2957.25
$ws.Range("N100").Value = -4904.5
# Row 138
$ws.Range("H138").Value = 47619.625
$ws.Range("I138").Value = 0
$ws.Range("J138").Value = 47619.625
$ws.Range("K138").Value = 0
$ws.Range("L138").Value = 47619.625
$ws.Range("N138").Value = -57899.625

# ---- Sheet: WVR ----
$ws = $wb.Worksheets.Item("WVR")
# Row 81
$ws.Range("H81").Value = 8104.8486
$ws.Range("I81").Value = 1058.1
$ws.Range("J81").Value = 18946
$ws.Range("K81").Value = 2116.2
$ws.Range("L81").Value = 37892
$ws.Range("M81").Value = -1055.2
$ws.Range("N81").Value = -40014
# Row 84
$ws.Range("H84").Value = 8104.8486
$ws.Range("I84").Value = 1058.1
$ws.Range("J84").Value = 18946
$ws.Range("K84").Value = 10581
$ws.Range("L84").Value = 189460
$ws.Range("M84").Value = -5277
$ws.Range("N84").Value = -200068
# Row 138
$ws.Range("H138").Value = 54260
$ws.Range("I138").Value = 0
$ws.Range("J138").Value = 54260
$ws.Range("K138").Value = 0
$ws.Range("L138").Value = 54260
$ws.Range("N138").Value = -64540
